# Insert a new weekly data row at row 9 (shifts old rows 9..95 down to 10..96),
# then fill the newly-inserted row 9 with this week's record -- identical to
# the (now shifted) row 10 except for the date (column D) and volume (column J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9 through 95 down by one row.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Cells.Item(9, 1).Value = 9
$ws.Cells.Item(9, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(9, 3).Value = "Metropolitana"
$ws.Cells.Item(9, 4).Value = 44685
$ws.Cells.Item(9, 5).Value = 13
$ws.Cells.Item(9, 6).Value = 100112005
$ws.Cells.Item(9, 7).Value = "Puerro"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 160
$ws.Cells.Item(9, 11).Value = 8000
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 8000
$ws.Cells.Item(9, 14).Value = "$/paquete 20 unidades"
$ws.Cells.Item(9, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(9, 16).Value = 400
$ws.Cells.Item(9, 17).Value = 20
$ws.Cells.Item(9, 18).Value = "Hortaliza"
